$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs1"
$ws.Range("C2").Value = "Itga3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.169512999999999
$ws.Range("N2").Value = 18.508539
$ws.Range("O2").Value = 0.5207942167525852
$ws.Range("P2").Value = 0.5207942167525853
$ws.Range("Q2").Value = 923.3068647221927
$ws.Range("R2").Value = 8309.761782499736
$ws.Range("S2").Value = 0.3083810457001705
$ws.Range("T2").Value = 0.3083810457001706

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs1"
$ws.Range("C3").Value = "Itga3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("N3").Value = 0.204412
$ws.Range("O3").Value = 0.005751755307905689
$ws.Range("P3").Value = 0.00575175530790569
$ws.Range("Q3").Value = 10.19718535491066
$ws.Range("R3").Value = 91.77466819419598
$ws.Range("S3").Value = 0.003405821837891324
$ws.Range("T3").Value = 0.003405821837891325

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Thbs1"
$ws.Range("C4").Value = "Itga3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.608704333333333
$ws.Range("N4").Value = 16.826113
$ws.Range("O4").Value = 0.473454027939509
$ws.Range("P4").Value = 0.4734540279395091
$ws.Range("Q4").Value = 839.3782804515974
$ws.Range("R4").Value = 7554.404524064378
$ws.Range("S4").Value = 0.2803492118966945
$ws.Range("T4").Value = 0.2803492118966945

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Thbs1"
$ws.Range("C5").Value = "Itga3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.169512999999999
$ws.Range("N5").Value = 18.508539
$ws.Range("O5").Value = 0.5207942167525852
$ws.Range("P5").Value = 0.5207942167525853
$ws.Range("Q5").Value = 325.3438676746206
$ws.Range("R5").Value = 2928.094809071586
$ws.Range("S5").Value = 0.1086636371493078
$ws.Range("T5").Value = 0.1086636371493079

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Thbs1"
$ws.Range("C6").Value = "Itga3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("N6").Value = 0.204412
$ws.Range("O6").Value = 0.005751755307905689
$ws.Range("P6").Value = 0.00575175530790569
$ws.Range("Q6").Value = 3.593162630454222
$ws.Range("R6").Value = 32.338463674088
$ws.Range("S6").Value = 0.001200102903690254
$ws.Range("T6").Value = 0.001200102903690254

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Thbs1"
$ws.Range("C7").Value = "Itga3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.608704333333333
$ws.Range("N7").Value = 16.826113
$ws.Range("O7").Value = 0.473454027939509
$ws.Range("P7").Value = 0.4734540279395091
$ws.Range("Q7").Value = 295.7701135324735
$ws.Range("R7").Value = 2661.931021792262
$ws.Range("S7").Value = 0.09878611367786791
$ws.Range("T7").Value = 0.09878611367786792

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Thbs1"
$ws.Range("C8").Value = "Itga3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.169512999999999
$ws.Range("N8").Value = 18.508539
$ws.Range("O8").Value = 0.5207942167525852
$ws.Range("P8").Value = 0.5207942167525853
$ws.Range("Q8").Value = 310.6308192417336
$ws.Range("R8").Value = 2795.677373175603
$ws.Range("S8").Value = 0.1037495339031069
$ws.Range("T8").Value = 0.1037495339031069

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Thbs1"
$ws.Range("C9").Value = "Itga3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("N9").Value = 0.204412
$ws.Range("O9").Value = 0.005751755307905689
$ws.Range("P9").Value = 0.00575175530790569
$ws.Range("Q9").Value = 3.430668786058222
$ws.Range("R9").Value = 30.876019074524
$ws.Range("S9").Value = 0.00114583056632411
$ws.Range("T9").Value = 0.001145830566324111

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Thbs1"
$ws.Range("C10").Value = "Itga3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.608704333333333
$ws.Range("N10").Value = 16.826113
$ws.Range("O10").Value = 0.473454027939509
$ws.Range("P10").Value = 0.4734540279395091
$ws.Range("Q10").Value = 282.3944810470445
$ws.Range("R10").Value = 2541.550329423401
$ws.Range("S10").Value = 0.09431870236494666
$ws.Range("T10").Value = 0.09431870236494669

